$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the Spanish prepositions ("de"/"del"/"el"/"la") inside a
#     handful of place names so they read "De"/"Del"/"El"/"La" ---
$ws.Range("A15").Value = "Ciudad De México"
$ws.Range("A25").Value = "Estado De México"
$ws.Range("B27").Value = "San Felipe Del Progreso"
$ws.Range("B36").Value = "Acapulco De Juárez"
$ws.Range("B45").Value = "Acatlán De Juárez"
$ws.Range("B46").Value = "Atotonilco El Alto"
$ws.Range("B50").Value = "La Manzanilla De La Paz"
$ws.Range("B54").Value = "San Cristóbal De La Barranca"
$ws.Range("B56").Value = "Tepatitlán De Morelos"
$ws.Range("B79").Value = "Chalcatongo De Hidalgo"
$ws.Range("B80").Value = "Constancia Del Rosario"
$ws.Range("B85").Value = "San Miguel El Grande"
$ws.Range("B93").Value = "Tetela De Ocampo"
$ws.Range("B95").Value = "Cadereyta De Montes"

# --- Drop the trailing footnote/source blocks that lived past the real
#     data (rows 115-119 and 476-480); clearing their contents removes
#     them from the used range and shrinks the sheet dimension down to
#     the real data at A1:D113 ---
$ws.Range("A115:D119").ClearContents()
$ws.Range("A476:D480").ClearContents()
